$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "M2"
$ws.Range("B2").Value = "Cd200r1"
$ws.Range("C2").Value = "Cd200"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.93046766666667
$ws.Range("H2").Value = 35.791403
$ws.Range("I2").Value = 0.9935182622400663
$ws.Range("J2").Value = 0.9935182622400662
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 65.569613
$ws.Range("N2").Value = 196.708839
$ws.Range("O2").Value = 0.6815983651189281
$ws.Range("P2").Value = 0.681598365118928
$ws.Range("Q2").Value = 782.2761478123464
$ws.Range("R2").Value = 7040.485330311118
$ws.Range("S2").Value = 0.6771804232586277
$ws.Range("T2").Value = 0.6771804232586275
$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "Cd200r1"
$ws.Range("C3").Value = "Cd200"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.93046766666667
$ws.Range("H3").Value = 35.791403
$ws.Range("I3").Value = 0.9935182622400663
$ws.Range("J3").Value = 0.9935182622400662
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.561623333333333
$ws.Range("N3").Value = 13.68487
$ws.Range("O3").Value = 0.04741823024467683
$ws.Range("P3").Value = 0.04741823024467683
$ws.Range("Q3").Value = 54.42229968584556
$ws.Range("R3").Value = 489.8006971726101
$ws.Range("S3").Value = 0.04711087771119068
$ws.Range("T3").Value = 0.04711087771119067
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Cd200r1"
$ws.Range("C4").Value = "Cd200"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.93046766666667
$ws.Range("H4").Value = 35.791403
$ws.Range("I4").Value = 0.9935182622400663
$ws.Range("J4").Value = 0.9935182622400662
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.611825
$ws.Range("N4").Value = 1.835475
$ws.Range("O4").Value = 0.006359941757455365
$ws.Range("P4").Value = 0.006359941757455365
$ws.Range("Q4").Value = 7.299358380158333
$ws.Range("R4").Value = 65.694225421425
$ws.Range("S4").Value = 0.006318718282815087
$ws.Range("T4").Value = 0.006318718282815087
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Cd200r1"
$ws.Range("C5").Value = "Cd200"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.93046766666667
$ws.Range("H5").Value = 35.791403
$ws.Range("I5").Value = 0.9935182622400663
$ws.Range("J5").Value = 0.9935182622400662
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 25.456719
$ws.Range("N5").Value = 76.37015699999999
$ws.Range("O5").Value = 0.2646234628789398
$ws.Range("P5").Value = 0.2646234628789398
$ws.Range("Q5").Value = 303.7105629289189
$ws.Range("R5").Value = 2733.395066360271
$ws.Range("S5").Value = 0.262908242987433
$ws.Range("T5").Value = 0.262908242987433
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Cd200r1"
$ws.Range("C6").Value = "Cd200"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.07783466666666666
$ws.Range("H6").Value = 0.233504
$ws.Range("I6").Value = 0.006481737759933703
$ws.Range("J6").Value = 0.006481737759933702
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 65.569613
$ws.Range("N6").Value = 196.708839
$ws.Range("O6").Value = 0.6815983651189281
$ws.Range("P6").Value = 0.681598365118928
$ws.Range("Q6").Value = 5.103588971317333
$ws.Range("R6").Value = 45.932300741856
$ws.Range("S6").Value = 0.004417941860300436
$ws.Range("T6").Value = 0.004417941860300434
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Cd200r1"
$ws.Range("C7").Value = "Cd200"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.07783466666666666
$ws.Range("H7").Value = 0.233504
$ws.Range("I7").Value = 0.006481737759933703
$ws.Range("J7").Value = 0.006481737759933702
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.561623333333333
$ws.Range("N7").Value = 13.68487
$ws.Range("O7").Value = 0.04741823024467683
$ws.Range("P7").Value = 0.04741823024467683
$ws.Range("Q7").Value = 0.3550524316088889
$ws.Range("R7").Value = 3.19547188448
$ws.Range("S7").Value = 0.0003073525334861522
$ws.Range("T7").Value = 0.0003073525334861521
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cd200r1"
$ws.Range("C8").Value = "Cd200"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.07783466666666666
$ws.Range("H8").Value = 0.233504
$ws.Range("I8").Value = 0.006481737759933703
$ws.Range("J8").Value = 0.006481737759933702
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.611825
$ws.Range("N8").Value = 1.835475
$ws.Range("O8").Value = 0.006359941757455365
$ws.Range("P8").Value = 0.006359941757455365
$ws.Range("Q8").Value = 0.04762119493333333
$ws.Range("R8").Value = 0.4285907544
$ws.Range("S8").Value = [double]"4.122347464027756e-05"
$ws.Range("T8").Value = [double]"4.122347464027755e-05"
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cd200r1"
$ws.Range("C9").Value = "Cd200"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.07783466666666666
$ws.Range("H9").Value = 0.233504
$ws.Range("I9").Value = 0.006481737759933703
$ws.Range("J9").Value = 0.006481737759933702
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 25.456719
$ws.Range("N9").Value = 76.37015699999999
$ws.Range("O9").Value = 0.2646234628789398
$ws.Range("P9").Value = 0.2646234628789398
$ws.Range("Q9").Value = 1.981415237792
$ws.Range("R9").Value = 17.832737140128
$ws.Range("S9").Value = 0.001715219891506839
$ws.Range("T9").Value = 0.001715219891506839
